# Auto-generated script to apply the Sophia_Profits leve-profit data refresh
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each affected row, currentAveragePrice / currentAveragePriceNQ/HQ,
# LevePriceNQ/HQ and LeveProfitNQ/HQ columns (H-N) are refreshed with the
# latest market-board derived values produced by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4046.3809
$ws.Range("I15").Value = 4046.3809
$ws.Range("K15").Value = 12139.1427
$ws.Range("M15").Value = -11970.1427
$ws.Range("H74").Value = 5636.8887
$ws.Range("I74").Value = 5636.8887
$ws.Range("K74").Value = 5636.8887
$ws.Range("M74").Value = -4700.8887
$ws.Range("H77").Value = 5636.8887
$ws.Range("I77").Value = 5636.8887
$ws.Range("K77").Value = 28184.4435
$ws.Range("M77").Value = -23504.4435
$ws.Range("H92").Value = 897.5
$ws.Range("I92").Value = 993.2727
$ws.Range("K92").Value = 993.2727
$ws.Range("M92").Value = 254.7273
$ws.Range("H96").Value = 443.44446
$ws.Range("I96").Value = 248.71428
$ws.Range("J96").Value = 1125
$ws.Range("K96").Value = 746.14284
$ws.Range("L96").Value = 3375
$ws.Range("M96").Value = 626.85716
$ws.Range("N96").Value = -6121
$ws.Range("H97").Value = 1846
$ws.Range("J97").Value = 1846
$ws.Range("L97").Value = 5538
$ws.Range("N97").Value = -6530
$ws.Range("H106").Value = 5000
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H125").Value = 749.5
$ws.Range("I125").Value = 499.33334
$ws.Range("K125").Value = 4494.00006
$ws.Range("M125").Value = -2034.00006
$ws.Range("H129").Value = 1249.75
$ws.Range("I129").Value = 1249.75
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 3749.25
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 1250.75
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 3279
$ws.Range("I131").Value = 1695
$ws.Range("J131").Value = 3675
$ws.Range("K131").Value = 5085
$ws.Range("L131").Value = 11025
$ws.Range("M131").Value = -45
$ws.Range("N131").Value = -21105
$ws.Range("H138").Value = 3151
$ws.Range("J138").Value = 3449
$ws.Range("L138").Value = 10347
$ws.Range("N138").Value = -20627

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7262.696
$ws.Range("I32").Value = 7262.696
$ws.Range("K32").Value = 7262.696
$ws.Range("M32").Value = -6975.696
$ws.Range("H61").Value = 1713.4166
$ws.Range("I61").Value = 1728.2727
$ws.Range("K61").Value = 1728.2727
$ws.Range("M61").Value = -1516.2727
$ws.Range("H74").Value = 8243.4
$ws.Range("I74").Value = 8048.3335
$ws.Range("K74").Value = 8048.3335
$ws.Range("M74").Value = -7174.3335
$ws.Range("H77").Value = 8243.4
$ws.Range("I77").Value = 8048.3335
$ws.Range("K77").Value = 40241.6675
$ws.Range("M77").Value = -35873.6675
$ws.Range("H97").Value = 588.44446
$ws.Range("I97").Value = 393.8
$ws.Range("J97").Value = 831.75
$ws.Range("K97").Value = 393.8
$ws.Range("L97").Value = 831.75
$ws.Range("M97").Value = 102.2
$ws.Range("N97").Value = -1823.75
$ws.Range("H136").Value = 1713.4166
$ws.Range("I136").Value = 1728.2727
$ws.Range("K136").Value = 5184.8181
$ws.Range("M136").Value = -2634.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3131.7
$ws.Range("I94").Value = 2370.6667
$ws.Range("J94").Value = 4273.25
$ws.Range("K94").Value = 2370.6667
$ws.Range("L94").Value = 4273.25
$ws.Range("M94").Value = -1919.6667
$ws.Range("N94").Value = -5175.25
$ws.Range("H99").Value = 1805.9
$ws.Range("I99").Value = 1432.375
$ws.Range("K99").Value = 1432.375
$ws.Range("M99").Value = 65.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 1002.5
$ws.Range("J13").Value = 1002.5
$ws.Range("L13").Value = 1002.5
$ws.Range("N13").Value = -1280.5
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50
$ws.Range("H58").Value = 2037.7142
$ws.Range("I58").Value = 2215.8
$ws.Range("K58").Value = 2215.8
$ws.Range("M58").Value = -2012.8
$ws.Range("H134").Value = 6528.8184
$ws.Range("I134").Value = 6648.5
$ws.Range("K134").Value = 19945.5
$ws.Range("M134").Value = -17410.5
$ws.Range("H136").Value = 2037.7142
$ws.Range("I136").Value = 2215.8
$ws.Range("K136").Value = 6647.400000000001
$ws.Range("M136").Value = -4097.400000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1416.6666
$ws.Range("I13").Value = 250
$ws.Range("J13").Value = 2000
$ws.Range("K13").Value = 750
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = -582
$ws.Range("N13").Value = -6336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H70").Value = 2219.2
$ws.Range("I70").Value = 2219.2
$ws.Range("K70").Value = 2219.2
$ws.Range("M70").Value = -1949.2
$ws.Range("H73").Value = 2219.2
$ws.Range("I73").Value = 2219.2
$ws.Range("K73").Value = 2219.2
$ws.Range("M73").Value = -1283.2
$ws.Range("H102").Value = 5210
$ws.Range("I102").Value = 4474.4
$ws.Range("K102").Value = 4474.4
$ws.Range("M102").Value = -2852.4
$ws.Range("H107").Value = 2042.4286
$ws.Range("I107").Value = 2829.6
$ws.Range("K107").Value = 2829.6
$ws.Range("M107").Value = -909.5999999999999
$ws.Range("H113").Value = 2977.6667
$ws.Range("I113").Value = 3099.875
$ws.Range("K113").Value = 3099.875
$ws.Range("M113").Value = -929.875
$ws.Range("H123").Value = 22500
$ws.Range("J123").Value = 22500
$ws.Range("L123").Value = 22500
$ws.Range("N123").Value = -27400
$ws.Range("H132").Value = 4004.375
$ws.Range("I132").Value = 3672.8333
$ws.Range("K132").Value = 11018.4999
$ws.Range("M132").Value = -8488.499899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3424.2856
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("N46").Value = -10376
$ws.Range("H61").Value = 994.5
$ws.Range("I61").Value = 999
$ws.Range("J61").Value = 990
$ws.Range("K61").Value = 999
$ws.Range("L61").Value = 990
$ws.Range("M61").Value = -797
$ws.Range("N61").Value = -1394
$ws.Range("H93").Value = 3817.5
$ws.Range("I93").Value = 3801
$ws.Range("J93").Value = 3834
$ws.Range("K93").Value = 3801
$ws.Range("L93").Value = 3834
$ws.Range("M93").Value = -2553
$ws.Range("N93").Value = -6330
$ws.Range("H100").Value = 4296.5
$ws.Range("I100").Value = 4815.8
$ws.Range("J100").Value = 1700
$ws.Range("K100").Value = 4815.8
$ws.Range("L100").Value = 1700
$ws.Range("M100").Value = -4274.8
$ws.Range("N100").Value = -2782
$ws.Range("H113").Value = 994.5
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 990
$ws.Range("K113").Value = 999
$ws.Range("L113").Value = 990
$ws.Range("M113").Value = 1171
$ws.Range("N113").Value = -5330

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1604.8235
$ws.Range("I100").Value = 1152.909
$ws.Range("K100").Value = 2305.818
$ws.Range("M100").Value = -1764.818
$ws.Range("H113").Value = 916.75
$ws.Range("I113").Value = 1187.625
$ws.Range("J113").Value = 375
$ws.Range("K113").Value = 3562.875
$ws.Range("L113").Value = 1125
$ws.Range("M113").Value = -1392.875
$ws.Range("N113").Value = -5465
$ws.Range("H122").Value = 4003.2778
$ws.Range("I122").Value = 3940.875
$ws.Range("K122").Value = 11822.625
$ws.Range("M122").Value = -9372.625
$ws.Range("H126").Value = 1264.9656
$ws.Range("I126").Value = 1046.4445
$ws.Range("K126").Value = 3139.3335
$ws.Range("M126").Value = -669.3335000000002

